$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove "/RME" from the steel description line in B2
$cell = $ws.Range("B2")
$cell.Value2 = $cell.Value2 -replace "/RME/H:1", "/H:1"

# Wrap text on B2 and resize row to fit the multi-line content
$cell.WrapText = $true
$ws.Rows.Item(2).RowHeight = 365

# Update selection to match the author's final selection (B2:B13, active cell B13)
$ws.Range("B2:B13").Select()
